$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3, 12, 3, 7, 6, 8, 3, 7, 6, 10, 7, 9, 6, 9, 5, 6, 8, 6, 12, 5, 4, 7, 6, 7, 11, 3, 6, 5, 10, 4, 4, 6, 7, 5, 3, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
